# !diary.xlsx - diary entry update for 11.11.18
# -Preparation of OutputOptions -NextStep Implementation of OutputOptions
#
# Changes applied:
#   1. Remarks (E21) for the 11.11.18 entry get the real diary text instead
#      of the "-" placeholder, wrapped over three lines.
#   2. End time (C21) moves from 21:00 to 20:00, which shifts the computed
#      duration (D21) and the grand total (D26) accordingly (both are
#      formulas, so they recalculate automatically).
#   3. Row 21 is made taller to fit the now-multiline remarks, and the
#      remarks cell gets word-wrap turned on (matching the formatting
#      already used on the other multi-line remark cells in the sheet).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("summary")

# New end time for 11.11.18: 20:00 instead of 21:00.
$ws.Range("C21").Value = 0.83333333333333337

# New remarks text (three lines).
$ws.Range("E21").Value = "-Added timer`n-Added Options for Idling after frame calculation`n-Prepared OutputOptions"

# Match the wrap-text formatting used by the other multi-line remark cells.
$ws.Range("E21").WrapText = $true

# Grow the row so the wrapped text is fully visible.
$ws.Rows("21").RowHeight = 45

Write-Output "Updated C21, E21 and row 21 height on sheet 'summary'."
